# Split the reports details in two different lists (second part)
#
# 1. Rename the existing (only) sheet from "Report C - Salida" to "flagged IM's".
# 2. Duplicate it (same columns/format/header row/empty trailing rows) to create
#    a second sheet named "undetermined", placed right after the first one.
# 3. Restore the first sheet ("flagged IM's") as the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "flagged IM's"

# Copy the sheet (with all its content/formatting) to just after itself.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "undetermined"

# Keep the original sheet as the selected/active one (Copy leaves the new
# sheet active).
$ws1.Activate()
$ws1.Select()
